# Regenerate orders with updated distance/size codes.
# The experiment's distance and size condition labels are renamed:
#   D80 -> D86, D64 -> D69, D51 -> D55, and the "S30" size token -> "S31"
# These tokens appear (as substrings) inside many shared text values across
# the sheet (Condition, Filename_Left, Filename_Right, Distance, Size
# columns), so we perform the substitutions across the whole used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange

$usedRange.Replace("D80", "D86", -4143, 2, $false, $false)
$usedRange.Replace("D64", "D69", -4143, 2, $false, $false)
$usedRange.Replace("D51", "D55", -4143, 2, $false, $false)
$usedRange.Replace("S30", "S31", -4143, 2, $false, $false)
